$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("B1").Value = "California"
$ws.Range("C1").Value = 45272
$ws.Range("C1").NumberFormat = "mm-dd-yy"

$ws.Range("A17").Value = "For petroleum, we use the weighted average expected capacity factor."
$ws.Range("A18").Value = "In reality, petroleum is primarily used in fuel-switching applications in New England"
$ws.Range("A19").Value = "(where there can be specific, unusual hours of NG shortage due to pipeline congestion)"
$ws.Range("A20").Value = "and in places like Hawaii, without much NG access, where it is easier to transport petroleum"
$ws.Range("A21").Value = "by ship due to its higher volumetric energy density."
$ws.Range("A22").Value = "We can't capture some of these unusual hours and sub-regional dynamics directly, so"
$ws.Range("A23").Value = "we use guaranteed dispatch to align modeled petroleum dispatch with real-world experience."
$ws.Range("A24:A26").ClearContents()
$ws3 = $wb.Worksheets.Item("BGDPbES")
$ws3.Range("B10:AK10").Value = 1
$ws3.Range("B11:AK11").Value = 1
